$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.138.16'

$ws.Range("D3").Value = '1.903.45'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.39'
$ws.Range("E5").Value = '  +2.91%  '

$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.72'
$ws.Range("E8").Value = '  +2.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.354'
$ws.Range("E9").Value = '  +2.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.33'
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("E11").Value = '  +5.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0980'
$ws.Range("E12").Value = '  -0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.26'
$ws.Range("E13").Value = '  +5.37%  '

$ws.Range("D14").Value = '2.181.84'
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.734'
$ws.Range("E15").Value = '  +4.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.01'
$ws.Range("E16").Value = '  +4.39%  '

$ws.Range("D17").Value = '1.916.37'
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").Value = '35.160.01'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.73'
$ws.Range("E19").Value = '  +2.48%  '

$ws.Range("D20").Value = '0.0₃0840'
$ws.Range("E20").Value = '  +2.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '243.28'
$ws.Range("E21").Value = '  +1.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.05'
$ws.Range("E22").Value = '  +2.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.04'
$ws.Range("E23").Value = '  +5.41%  '

$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("E25").Value = '  +4.78%  '

$ws.Range("E26").Value = '  -1.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.30'
$ws.Range("E27").Value = '  +0.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.60'
$ws.Range("E28").Value = '  -0.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.52'
$ws.Range("E29").Value = '  +1.13%  '

$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("D31").Value = '4.128.24'
$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("E32").Value = '  +10.92%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.34'
$ws.Range("E33").Value = '  +4.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0597'
$ws.Range("E34").Value = '  +5.48%  '

$ws.Range("E35").Value = '  +8.37%  '

$ws.Range("E36").Value = '  +4.07%  '

$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.850'
$ws.Range("E38").Value = '  -6.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.01'
$ws.Range("E39").Value = '  -0.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.30'
$ws.Range("E40").Value = '  +5.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.79'
$ws.Range("E41").Value = '  +6.99%  '

$ws.Range("E42").Value = '  +4.37%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0664'
$ws.Range("E43").Value = '  +1.56%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.11'
$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").Value = '1.303.83'
$ws.Range("E46").Value = '  -3.22%  '

$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.60'
$ws.Range("E49").Value = '  +2.06%  '

$ws.Range("B50").Value = 'Gas'
$ws.Range("C50").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.03'
$ws.Range("E50").Value = '  -0.70%  '

$ws.Range("E51").Value = '  +7.13%  '
